# "Heuristique mise en forme"
# - Rename the "Heuristique v2" column header (H2) to "Heuristique v3"
# - For every data row that has both a "Notre Heuristique" (col G) value and a
#   "Heuristique v2" (col H) value, drop the old G value, move the H value
#   into G, and clear out the H cell entirely.
# - Update the active selection to H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column header H2: "Heuristique v2" -> "Heuristique v3"
$ws.Cells.Item(2, 8).Value = "Heuristique v3"

# Rows where both G (col 7) and H (col 8) had values: move H's value into G,
# then clear H.
$rowsToMerge = @{
    4  = 798
    5  = 810
    6  = 432
    7  = 816
    8  = 540
    9  = 492
    10 = 570
    11 = 414
    12 = 750
    13 = 318
    14 = 590
    17 = 367
    18 = 385
    21 = 338
    22 = 257
    23 = 311
    24 = 256
    25 = 264
    26 = 466
    27 = 260
    28 = 178
    29 = 720
}

foreach ($row in $rowsToMerge.Keys) {
    $ws.Cells.Item($row, 7).Value = $rowsToMerge[$row]
    $ws.Cells.Item($row, 8).ClearContents()
}

# Update selection to H2 (matches the saved sheetView selection in the diff)
$ws.Range("H2").Select()
